# Insert a new weekly price record for Cilantro - Agrícola del Norte S.A. de Arica
# This shifts the existing rows 59..93 down to 60..94 and populates the
# newly inserted row 59 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59 (existing row 59 and below shift down by one)
$ws.Rows(59).Insert()

# Populate the new row 59 with the new weekly record
$ws.Cells.Item(59, 1).Value = 1
$ws.Cells.Item(59, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(59, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(59, 4).Value = 44879
$ws.Cells.Item(59, 5).Value = 15
$ws.Cells.Item(59, 6).Value = 100112040
$ws.Cells.Item(59, 7).Value = "Cilantro"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 300
$ws.Cells.Item(59, 11).Value = 1300
$ws.Cells.Item(59, 12).Value = 1500
$ws.Cells.Item(59, 13).Value = 1400
$ws.Cells.Item(59, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(59, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(59, 16).Value = 700
$ws.Cells.Item(59, 17).Value = 2
$ws.Cells.Item(59, 18).Value = "Hortaliza"
